# Update TPM-derived NATMI LR-pair metrics (Mdk-Alk) to new normalization values.
# Only the numeric result columns (G..T) change; identifiers in A..F are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.988074333333333  # G2
$ws.Cells.Item(2, 8).Value = 5.964223  # H2
$ws.Cells.Item(2, 9).Value = 0.01657769708907969  # I2
$ws.Cells.Item(2, 10).Value = 0.01657769708907968  # J2
$ws.Cells.Item(2, 15).Value = 0.3717075934090293  # O2
$ws.Cells.Item(2, 16).Value = 0.3717075934090293  # P2
$ws.Cells.Item(2, 17).Value = 0.1000114047212222  # Q2
$ws.Cells.Item(2, 18).Value = 0.9001026424909999  # R2
$ws.Cells.Item(2, 19).Value = 0.006162055889245681  # S2
$ws.Cells.Item(2, 20).Value = 0.00616205588924568  # T2

# Row 3
$ws.Cells.Item(3, 7).Value = 1.988074333333333  # G3
$ws.Cells.Item(3, 8).Value = 5.964223  # H3
$ws.Cells.Item(3, 9).Value = 0.01657769708907969  # I3
$ws.Cells.Item(3, 10).Value = 0.01657769708907968  # J3
$ws.Cells.Item(3, 11).Value = 3  # K3
$ws.Cells.Item(3, 12).Value = 1  # L3
$ws.Cells.Item(3, 13).Value = 0.08503100000000001  # M3
$ws.Cells.Item(3, 14).Value = 0.255093  # N3
$ws.Cells.Item(3, 15).Value = 0.6282924065909707  # O3
$ws.Cells.Item(3, 16).Value = 0.6282924065909707  # P3
$ws.Cells.Item(3, 17).Value = 0.1690479486376667  # Q3
$ws.Cells.Item(3, 18).Value = 1.521431537739  # R3
$ws.Cells.Item(3, 19).Value = 0.01041564119983401  # S3
$ws.Cells.Item(3, 20).Value = 0.01041564119983401  # T3

# Row 4
$ws.Cells.Item(4, 9).Value = 0.7746030815641455  # I4
$ws.Cells.Item(4, 10).Value = 0.7746030815641454  # J4
$ws.Cells.Item(4, 15).Value = 0.3717075934090293  # O4
$ws.Cells.Item(4, 16).Value = 0.3717075934090293  # P4
$ws.Cells.Item(4, 19).Value = 0.2879258472954266  # S4
$ws.Cells.Item(4, 20).Value = 0.2879258472954265  # T4

# Row 5
$ws.Cells.Item(5, 9).Value = 0.7746030815641455  # I5
$ws.Cells.Item(5, 10).Value = 0.7746030815641454  # J5
$ws.Cells.Item(5, 11).Value = 3  # K5
$ws.Cells.Item(5, 12).Value = 1  # L5
$ws.Cells.Item(5, 13).Value = 0.08503100000000001  # M5
$ws.Cells.Item(5, 14).Value = 0.255093  # N5
$ws.Cells.Item(5, 15).Value = 0.6282924065909707  # O5
$ws.Cells.Item(5, 16).Value = 0.6282924065909707  # P5
$ws.Cells.Item(5, 17).Value = 7.898869260501335  # Q5
$ws.Cells.Item(5, 18).Value = 71.089823344512  # R5
$ws.Cells.Item(5, 19).Value = 0.486677234268719  # S5
$ws.Cells.Item(5, 20).Value = 0.4866772342687189  # T5

# Row 6
$ws.Cells.Item(6, 7).Value = 23.741365  # G6
$ws.Cells.Item(6, 8).Value = 71.224095  # H6
$ws.Cells.Item(6, 9).Value = 0.1979690350870239  # I6
$ws.Cells.Item(6, 10).Value = 0.1979690350870239  # J6
$ws.Cells.Item(6, 15).Value = 0.3717075934090293  # O6
$ws.Cells.Item(6, 16).Value = 0.3717075934090293  # P6
$ws.Cells.Item(6, 17).Value = 1.194325193901667  # Q6
$ws.Cells.Item(6, 18).Value = 10.748926745115  # R6
$ws.Cells.Item(6, 19).Value = 0.07358659360170534  # S6
$ws.Cells.Item(6, 20).Value = 0.07358659360170534  # T6

# Row 7
$ws.Cells.Item(7, 7).Value = 23.741365  # G7
$ws.Cells.Item(7, 8).Value = 71.224095  # H7
$ws.Cells.Item(7, 9).Value = 0.1979690350870239  # I7
$ws.Cells.Item(7, 10).Value = 0.1979690350870239  # J7
$ws.Cells.Item(7, 11).Value = 3  # K7
$ws.Cells.Item(7, 12).Value = 1  # L7
$ws.Cells.Item(7, 13).Value = 0.08503100000000001  # M7
$ws.Cells.Item(7, 14).Value = 0.255093  # N7
$ws.Cells.Item(7, 15).Value = 0.6282924065909707  # O7
$ws.Cells.Item(7, 16).Value = 0.6282924065909707  # P7
$ws.Cells.Item(7, 17).Value = 2.018752007315  # Q7
$ws.Cells.Item(7, 18).Value = 18.168768065835  # R7
$ws.Cells.Item(7, 19).Value = 0.1243824414853186  # S7
$ws.Cells.Item(7, 20).Value = 0.1243824414853186  # T7

# Row 8
$ws.Cells.Item(8, 7).Value = 1.301204666666667  # G8
$ws.Cells.Item(8, 8).Value = 3.903614  # H8
$ws.Cells.Item(8, 9).Value = 0.01085018625975097  # I8
$ws.Cells.Item(8, 10).Value = 0.01085018625975097  # J8
$ws.Cells.Item(8, 15).Value = 0.3717075934090293  # O8
$ws.Cells.Item(8, 16).Value = 0.3717075934090293  # P8
$ws.Cells.Item(8, 17).Value = 0.06545796822644445  # Q8
$ws.Cells.Item(8, 18).Value = 0.589121714038  # R8
$ws.Cells.Item(8, 19).Value = 0.00403309662265175  # S8
$ws.Cells.Item(8, 20).Value = 0.00403309662265175  # T8

# Row 9
$ws.Cells.Item(9, 7).Value = 1.301204666666667  # G9
$ws.Cells.Item(9, 8).Value = 3.903614  # H9
$ws.Cells.Item(9, 9).Value = 0.01085018625975097  # I9
$ws.Cells.Item(9, 10).Value = 0.01085018625975097  # J9
$ws.Cells.Item(9, 11).Value = 3  # K9
$ws.Cells.Item(9, 12).Value = 1  # L9
$ws.Cells.Item(9, 13).Value = 0.08503100000000001  # M9
$ws.Cells.Item(9, 14).Value = 0.255093  # N9
$ws.Cells.Item(9, 15).Value = 0.6282924065909707  # O9
$ws.Cells.Item(9, 16).Value = 0.6282924065909707  # P9
$ws.Cells.Item(9, 17).Value = 0.1106427340113334  # Q9
$ws.Cells.Item(9, 18).Value = 0.995784606102  # R9
$ws.Cells.Item(9, 19).Value = 0.006817089637099222  # S9
$ws.Cells.Item(9, 20).Value = 0.006817089637099221  # T9
